# Add a new "2022-Q4" sheet (most recent quarter) in front of the existing
# "2022-Q3" sheet, and update the "总计" (totals) summary sheet with the new
# quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it
#    inherits the exact same column layout / header styling), then
#    rename it and place it directly before "2022-Q3".
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The fund-code / ratio columns (B..G) hold values that look numeric
# ("009693", "12.47", "0.6410", ...) but must stay text so leading
# zeros and fixed decimal places survive. Force the range to Text
# format before writing so Excel doesn't auto-coerce them to numbers.
$q4.Range("B2:G8").NumberFormat = "@"

# Row 2
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "009693"
$q4.Cells.Item(2,3).Value = "富国积极成长一年定期开放混合"
$q4.Cells.Item(2,4).Value = "12.47"
$q4.Cells.Item(2,5).Value = "98.05"
$q4.Cells.Item(2,6).Value = "5.14"
$q4.Cells.Item(2,7).Value = "0.6410"
$q4.Cells.Item(2,8).Value = 3

# Row 3
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "001985"
$q4.Cells.Item(3,3).Value = "富国低碳新经济混合A"
$q4.Cells.Item(3,4).Value = "18.93"
$q4.Cells.Item(3,5).Value = "93.83"
$q4.Cells.Item(3,6).Value = "2.50"
$q4.Cells.Item(3,7).Value = "0.4732"
$q4.Cells.Item(3,8).Value = 10

# Row 4
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "257020"
$q4.Cells.Item(4,3).Value = "国联安精选混合"
$q4.Cells.Item(4,4).Value = "9.67"
$q4.Cells.Item(4,5).Value = "92.82"
$q4.Cells.Item(4,6).Value = "4.18"
$q4.Cells.Item(4,7).Value = "0.4042"
$q4.Cells.Item(4,8).Value = 6

# Row 5
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "006864"
$q4.Cells.Item(5,3).Value = "国联安核心资产策略混合"
$q4.Cells.Item(5,4).Value = "4.67"
$q4.Cells.Item(5,5).Value = "91.48"
$q4.Cells.Item(5,6).Value = "5.19"
$q4.Cells.Item(5,7).Value = "0.2424"
$q4.Cells.Item(5,8).Value = 3

# Row 6
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "014325"
$q4.Cells.Item(6,3).Value = "国联安核心趋势一年持有混合A"
$q4.Cells.Item(6,4).Value = "3.78"
$q4.Cells.Item(6,5).Value = "90.40"
$q4.Cells.Item(6,6).Value = "3.24"
$q4.Cells.Item(6,7).Value = "0.1225"
$q4.Cells.Item(6,8).Value = 10

# Row 7
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "011306"
$q4.Cells.Item(7,3).Value = "富国低碳新经济混合C"
$q4.Cells.Item(7,4).Value = "1.87"
$q4.Cells.Item(7,5).Value = "93.83"
$q4.Cells.Item(7,6).Value = "2.50"
$q4.Cells.Item(7,7).Value = "0.0468"
$q4.Cells.Item(7,8).Value = 10

# Row 8
$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).Value = "014326"
$q4.Cells.Item(8,3).Value = "国联安核心趋势一年持有混合C"
$q4.Cells.Item(8,4).Value = "0.38"
$q4.Cells.Item(8,5).Value = "90.40"
$q4.Cells.Item(8,6).Value = "3.24"
$q4.Cells.Item(8,7).Value = "0.0123"
$q4.Cells.Item(8,8).Value = 10

# ------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q4 at
#    the top of the data (row 2), push the existing rows down one,
#    and append the 2022-Q1 row that rolls off the bottom.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 7
$total.Cells.Item(2,4).Value = 1.94

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 7
$total.Cells.Item(3,4).Value = 1.41

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 9
$total.Cells.Item(4,4).Value = 1.62

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2022-Q1"
$total.Cells.Item(5,3).Value = 7
$total.Cells.Item(5,4).Value = 1.31
